$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 42, shifting existing rows 42-74 down by one
# (and carrying their formatting down with them).
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new TODO item. The formatting
# (styles) for row 42 was already carried down by Insert() from the row that
# used to be there, so there's no need to touch .Style here.
$ws.Range("A42").Value2 = "Not done"
$ws.Range("B42").Value2 = "Check SunEye Obstructions and both Solar Pathfinder file imports"
$ws.Range("C42").Value2 = "Janine"

# Fix up the SUM formula in H17 to include the newly inserted row.
$ws.Range("H17").Formula = "=SUM(D17:D43)"

# Update the view: no frozen/top-left offset, selection on A43.
$ws.Application.ActiveWindow.ScrollRow = 1
[void]$ws.Range("A43").Select()
